{"js": "// Apply the set of text replacements described by the commit diff.\n// Each entry: exact text to search for (case sensitive) -> replacement text.\nconst replacements = [\n  {\n    find: \"Atmospheric pollution\",\n    replace: \"Atmospheric pollution I\"\n  },\n  {\n    find: \"Ativa\u00e7\u00e3o: 01/01/2022\",\n    replace: \"Ativa\u00e7\u00e3o: 01/01/2025\"\n  },\n  {\n    find: \"Enable the student to identify the main pollutants from the atmosphere, to evaluate their influence on climate change and make the management and monitoring of air quality in large Brazilian cities.\",\n    replace: \"Enable the student to identify the main pollutants in the atmosphere and their effects on man and the environment\"\n  },\n  {\n    find: \"7455355 - Robson da Silva Rocha\",\n    replace: \"4893449 - D\u00e9bora Souza Alvim\"\n  },\n  {\n    find: \"Caracter\u00edsticas e composi\u00e7\u00e3o da atmosfera. Origem, movimenta\u00e7\u00e3o e destino dos poluentes. Hist\u00f3rico da polui\u00e7\u00e3o do ar. Principais poluentes atmosf\u00e9ricos e padr\u00f5es da qualidade do ar.  Dispers\u00e3o de poluentes na atmosfera. Modelos matem\u00e1ticos do transporte de poluentes atmosf\u00e9ricos. Qualidade do ar no interior de edif\u00edcios. Controle de polui\u00e7\u00e3o de fontes fixas e m\u00f3veis. Otimiza\u00e7\u00e3o exergoecon\u00f4mica/ambiental.\",\n    replace: \"1)Caracter\u00edsticas e composi\u00e7\u00e3o da atmosfera. Origem, movimenta\u00e7\u00e3o e destino dos poluentes. Hist\u00f3rico da polui\u00e7\u00e3o do ar. Principais poluentes atmosf\u00e9ricos e padr\u00f5es da qualidade do ar.  O efeito estufa. Dispers\u00e3o de poluentes na atmosfera. Modelos matem\u00e1ticos do transporte de poluentes atmosf\u00e9ricos.  Controle da polui\u00e7\u00e3o atmosf\u00e9rica de fontes fixas e m\u00f3veis.\"\n  },\n  {\n    find: \"Characteristics and composition of the atmosphere. Origin, movement and fate of pollutants. History of air pollution. Major air pollutants and air quality standards. Dispersion of pollutants in the atmosphere. Mathematical models of transport of air pollutants. Air quality inside buildings. Control pollution of fixed and mobile sources. Exergoeconomic / environmental optimization\",\n    replace: \"1)Caracter\u00edsticas e composi\u00e7\u00e3o da atmosfera. Origem, movimenta\u00e7\u00e3o e destino dos poluentes. Hist\u00f3rico da polui\u00e7\u00e3o do ar. Principais poluentes atmosf\u00e9ricos e padr\u00f5es da qualidade do ar.  O efeito estufa. Dispers\u00e3o de poluentes na atmosfera. Modelos matem\u00e1ticos do transporte de poluentes atmosf\u00e9ricos.  Controle da polui\u00e7\u00e3o atmosf\u00e9rica de fontes fixas e m\u00f3veis.1)Characteristics and composition of the atmosphere. Origin, movement and fate of pollutants. History of air pollution. Major air pollutants and air quality standards. The greenhouse effect. Dispersion of pollutants in the atmosphere. Mathematical models of transport of air pollutants. Air Pollution control from fixed and mobile sources.\"\n  },\n  {\n    find: \"1) Caracteriza\u00e7\u00e3o da atmosfera e seus poluentes. 2) Padr\u00f5es da qualidade do ar. 3) Dispers\u00e3o de poluentes na atmosfera. 4) Modelagem matem\u00e1tica do transporte de poluentes.5) Qualidade do ar no interior de edif\u00edcios.6) Controle de polui\u00e7\u00e3o de fontes fixas e m\u00f3veis.7) Otimiza\u00e7\u00e3o exergoecon\u00f4mica/ambiental.\",\n    replace: \"Caracteriza\u00e7\u00e3o da atmosfera e seus poluentes. 2) Padr\u00f5es da qualidade do ar. 3) Dispers\u00e3o de poluentes na atmosfera. 4) O efeito estufa. 5) Modelagem matem\u00e1tica do transporte de poluentes. 6) Controle da polui\u00e7\u00e3o atmosf\u00e9rica de fontes fixas e m\u00f3veisA disciplina pode contar com viagens did\u00e1ticas para complementa\u00e7\u00e3o do conte\u00fado da disciplina.\"\n  },\n  {\n    find: \"Characterization of the atmosphere and its pollutants.Air quality standards.Dispersion of pollutants in the atmosphere.Mathematical modeling of pollutant transport.Air quality inside buildings.Control pollution of fixed and mobile sourcesExergoeconomic / environmental optimization.\",\n    replace: \"Characterization of the atmosphere and its pollutants. 2) Air quality standards. 3) Dispersion of pollutants in the atmosphere. 4) The greenhouse effect. 5) Mathematical modeling of pollutant transport. 6) Air Pollution control from fixed and mobile sources. The discipline may have didactic trips to complement the content of the discipline.\"\n  },\n  {\n    find: \"Bibliografia b\u00e1sica:GUNTER, F.; Introdu\u00e7\u00e3o aos problemas da polui\u00e7\u00e3o ambiental. 1 ed. S\u00e3o Paulo: Editora EPU, 2008.LENZI, E. F.; FAVERO, L.O.B. Introdu\u00e7\u00e3o \u00e0 qu\u00edmica da atmosfera  Ci\u00eancia, vida e sobreviv\u00eancia. 1\u00aa. ed. Rio de Janeiro: Editora LCT, 465p. 2009.SPIRO, T. G.; STIGLIANI, E. W. M. Qu\u00edmica ambiental. 2 ed. Sao Paulo: Pearson / Prentice Hall. 2008. 352p.VESILIND, P. A.; MORGAN, S. M., revis\u00e3o t\u00e9cnica Carlos Alberto de Moya Figueira Netto; Lineu Belico dos Reis. Introdu\u00e7\u00e3o \u00e0 Engenharia Ambiental. Tradu\u00e7\u00e3o da 2\u00aa edi\u00e7\u00e3o norte-americana. Editora Cengage Learning, S\u00e3o Paulo, 2015.Bibliografia complementar:JACOBSON, Mark Z. Atmospheric pollution: history, science, and regulation. Cambridge, Inglaterra: Cambridge University Press, c2002. xi, 399 p. Includes bibliographical references (p 355-370). ISBN 9780521010443.SEINFELD, J.H.;MANAHAN, S.E. Environmental chemistry. 9 th edition. Boca Raton, FL: CRC Press. 753p. 2010.PANDIS, S.N. Atmospheric Chemistry and Physics: From air pollution to climate change. John Wiley& Sons, 1998.SCHNELLE JR, Karl B; BROWN, Charles A. Air pollution control technology handbook. New York: CRC Press, 2001. 386 p. (Mechanical engineering handbook series). ISBN 9780849395888.VALLERO, Daniel A. Fundamentals of air pollution. 4 ed. Amsterdam: Elsevier, 2008. 942 p\",\n    replace: \"Bibliografia b\u00e1sica:Baird, C.; Cann, M. Qu\u00edmica Ambiental. Porto Alegre: Bookman, 4.ed., 2011. 844p.GUNTER, F.; Introdu\u00e7\u00e3o aos problemas da polui\u00e7\u00e3o ambiental. 1 ed. S\u00e3o Paulo: Editora EPU, 2008.LENZI, E. F.; FAVERO, L.O.B. Introdu\u00e7\u00e3o \u00e0 qu\u00edmica da atmosfera  Ci\u00eancia, vida e sobreviv\u00eancia. 1\u00aa. ed. Rio de Janeiro: Editora LCT, 465p. 2009.Rocha, Julio Cesar; Rosa, Andr\u00e9 Henrique; Cardoso, Arnaldo Alves. Introdu\u00e7\u00e3o \u00e0 qu\u00edmica ambiental. 2. ed. Porto Alegre: Bookman, 2009.  03Seinfeld, J.H. e Pandis, S.P. Atmospheric Chemistry and Physics: from air pollution to climate change. New York, USA: John Wiley & Sons Inc., 2006.SPIRO, T. G.; STIGLIANI, E. W. M. Qu\u00edmica ambiental. 2 ed. Sao Paulo: Pearson / Prentice Hall. 2008. 352p.\"\n  }\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + find.substring(0, 60));\n  }\n\n  // Replace only the first (and expected only) match to avoid touching\n  // unrelated occurrences of the same substring elsewhere in the document.\n  results.items[0].insertText(replace, \"Replace\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-FirstMatch($findText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $findText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $found = $find.Execute()\n    if ($found) {\n        $range.Text = $replaceText\n    } else {\n        throw \"No match found for: $findText\"\n    }\n}\n\nReplace-FirstMatch \"Atmospheric pollution\" \"Atmospheric pollution I\"\n\nReplace-FirstMatch \"Ativa\u00e7\u00e3o: 01/01/2022\" \"Ativa\u00e7\u00e3o: 01/01/2025\"\n\nReplace-FirstMatch \"Enable the student to identify the main pollutants from the atmosphere, to evaluate their influence on climate change and make the management and monitoring of air quality in large Brazilian cities.\" \"Enable the student to identify the main pollutants in the atmosphere and their effects on man and the environment\"\n\nReplace-FirstMatch \"7455355 - Robson da Silva Rocha\" \"4893449 - D\u00e9bora Souza Alvim\"\n\nReplace-FirstMatch \"Caracter\u00edsticas e composi\u00e7\u00e3o da atmosfera. Origem, movimenta\u00e7\u00e3o e destino dos poluentes. Hist\u00f3rico da polui\u00e7\u00e3o do ar. Principais poluentes atmosf\u00e9ricos e padr\u00f5es da qualidade do ar.  Dispers\u00e3o de poluentes na atmosfera. Modelos matem\u00e1ticos do transporte de poluentes atmosf\u00e9ricos. Qualidade do ar no interior de edif\u00edcios. Controle de polui\u00e7\u00e3o de fontes fixas e m\u00f3veis. Otimiza\u00e7\u00e3o exergoecon\u00f4mica/ambiental.\" \"1)Caracter\u00edsticas e composi\u00e7\u00e3o da atmosfera. Origem, movimenta\u00e7\u00e3o e destino dos poluentes. Hist\u00f3rico da polui\u00e7\u00e3o do ar. Principais poluentes atmosf\u00e9ricos e padr\u00f5es da qualidade do ar.  O efeito estufa. Dispers\u00e3o de poluentes na atmosfera. Modelos matem\u00e1ticos do transporte de poluentes atmosf\u00e9ricos.  Controle da polui\u00e7\u00e3o atmosf\u00e9rica de fontes fixas e m\u00f3veis.\"\n\nReplace-FirstMatch \"Characteristics and composition of the atmosphere. Origin, movement and fate of pollutants. History of air pollution. Major air pollutants and air quality standards. Dispersion of pollutants in the atmosphere. Mathematical models of transport of air pollutants. Air quality inside buildings. Control pollution of fixed and mobile sources. Exergoeconomic / environmental optimization\" \"1)Caracter\u00edsticas e composi\u00e7\u00e3o da atmosfera. Origem, movimenta\u00e7\u00e3o e destino dos poluentes. Hist\u00f3rico da polui\u00e7\u00e3o do ar. Principais poluentes atmosf\u00e9ricos e padr\u00f5es da qualidade do ar.  O efeito estufa. Dispers\u00e3o de poluentes na atmosfera. Modelos matem\u00e1ticos do transporte de poluentes atmosf\u00e9ricos.  Controle da polui\u00e7\u00e3o atmosf\u00e9rica de fontes fixas e m\u00f3veis.1)Characteristics and composition of the atmosphere. Origin, movement and fate of pollutants. History of air pollution. Major air pollutants and air quality standards. The greenhouse effect. Dispersion of pollutants in the atmosphere. Mathematical models of transport of air pollutants. Air Pollution control from fixed and mobile sources.\"\n\nReplace-FirstMatch \"1) Caracteriza\u00e7\u00e3o da atmosfera e seus poluentes. 2) Padr\u00f5es da qualidade do ar. 3) Dispers\u00e3o de poluentes na atmosfera. 4) Modelagem matem\u00e1tica do transporte de poluentes.5) Qualidade do ar no interior de edif\u00edcios.6) Controle de polui\u00e7\u00e3o de fontes fixas e m\u00f3veis.7) Otimiza\u00e7\u00e3o exergoecon\u00f4mica/ambiental.\" \"Caracteriza\u00e7\u00e3o da atmosfera e seus poluentes. 2) Padr\u00f5es da qualidade do ar. 3) Dispers\u00e3o de poluentes na atmosfera. 4) O efeito estufa. 5) Modelagem matem\u00e1tica do transporte de poluentes. 6) Controle da polui\u00e7\u00e3o atmosf\u00e9rica de fontes fixas e m\u00f3veisA disciplina pode contar com viagens did\u00e1ticas para complementa\u00e7\u00e3o do conte\u00fado da disciplina.\"\n\nReplace-FirstMatch \"Characterization of the atmosphere and its pollutants.Air quality standards.Dispersion of pollutants in the atmosphere.Mathematical modeling of pollutant transport.Air quality inside buildings.Control pollution of fixed and mobile sourcesExergoeconomic / environmental optimization.\" \"Characterization of the atmosphere and its pollutants. 2) Air quality standards. 3) Dispersion of pollutants in the atmosphere. 4) The greenhouse effect. 5) Mathematical modeling of pollutant transport. 6) Air Pollution control from fixed and mobile sources. The discipline may have didactic trips to complement the content of the discipline.\"\n\nReplace-FirstMatch \"Bibliografia b\u00e1sica:GUNTER, F.; Introdu\u00e7\u00e3o aos problemas da polui\u00e7\u00e3o ambiental. 1 ed. S\u00e3o Paulo: Editora EPU, 2008.LENZI, E. F.; FAVERO, L.O.B. Introdu\u00e7\u00e3o \u00e0 qu\u00edmica da atmosfera  Ci\u00eancia, vida e sobreviv\u00eancia. 1\u00aa. ed. Rio de Janeiro: Editora LCT, 465p. 2009.SPIRO, T. G.; STIGLIANI, E. W. M. Qu\u00edmica ambiental. 2 ed. Sao Paulo: Pearson / Prentice Hall. 2008. 352p.VESILIND, P. A.; MORGAN, S. M., revis\u00e3o t\u00e9cnica Carlos Alberto de Moya Figueira Netto; Lineu Belico dos Reis. Introdu\u00e7\u00e3o \u00e0 Engenharia Ambiental. Tradu\u00e7\u00e3o da 2\u00aa edi\u00e7\u00e3o norte-americana. Editora Cengage Learning, S\u00e3o Paulo, 2015.Bibliografia complementar:JACOBSON, Mark Z. Atmospheric pollution: history, science, and regulation. Cambridge, Inglaterra: Cambridge University Press, c2002. xi, 399 p. Includes bibliographical references (p 355-370). ISBN 9780521010443.SEINFELD, J.H.;MANAHAN, S.E. Environmental chemistry. 9 th edition. Boca Raton, FL: CRC Press. 753p. 2010.PANDIS, S.N. Atmospheric Chemistry and Physics: From air pollution to climate change. John Wiley& Sons, 1998.SCHNELLE JR, Karl B; BROWN, Charles A. Air pollution control technology handbook. New York: CRC Press, 2001. 386 p. (Mechanical engineering handbook series). ISBN 9780849395888.VALLERO, Daniel A. Fundamentals of air pollution. 4 ed. Amsterdam: Elsevier, 2008. 942 p\" \"Bibliografia b\u00e1sica:Baird, C.; Cann, M. Qu\u00edmica Ambiental. Porto Alegre: Bookman, 4.ed., 2011. 844p.GUNTER, F.; Introdu\u00e7\u00e3o aos problemas da polui\u00e7\u00e3o ambiental. 1 ed. S\u00e3o Paulo: Editora EPU, 2008.LENZI, E. F.; FAVERO, L.O.B. Introdu\u00e7\u00e3o \u00e0 qu\u00edmica da atmosfera  Ci\u00eancia, vida e sobreviv\u00eancia. 1\u00aa. ed. Rio de Janeiro: Editora LCT, 465p. 2009.Rocha, Julio Cesar; Rosa, Andr\u00e9 Henrique; Cardoso, Arnaldo Alves. Introdu\u00e7\u00e3o \u00e0 qu\u00edmica ambiental. 2. ed. Porto Alegre: Bookman, 2009.  03Seinfeld, J.H. e Pandis, S.P. Atmospheric Chemistry and Physics: from air pollution to climate change. New York, USA: John Wiley & Sons Inc., 2006.SPIRO, T. G.; STIGLIANI, E. W. M. Qu\u00edmica ambiental. 2 ed. Sao Paulo: Pearson / Prentice Hall. 2008. 352p.\"\n"}
